# rmm corrigida para 2022
# Adds the 2022 row (A28:E28) to the RMM table, with RMM_C (E28) = 57.7 and
# MM_declaradas/NV/RMM_NC (B28:D28) left blank, matching the border/format
# treatment used for the last row of the table (row 27 loses its bottom
# border so it no longer looks like the final row, and row 28 gets the
# closing border instead).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row 27 (2021) is no longer the last row of the table: drop the
#    bottom border on B27:D27 so only the new row 28 carries it.
#    (A27 keeps its border untouched, E27's border/format is unchanged.)
# ---------------------------------------------------------------------
$ws.Range("B27:D27").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none

# ---------------------------------------------------------------------
# 2) Build row 28 (2022).
# ---------------------------------------------------------------------

# A28: same look as the header-style "boxed" cells (font, alignment,
# full surrounding border) used elsewhere in column A.
$ws.Range("A1").Copy()
$ws.Range("A28").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A28").Value = 2022

# B28: blank cell, just gets a left/top/bottom medium border (no right
# border, since C28 supplies the shared line).
$rngB28 = $ws.Range("B28")
$rngB28.Borders.Weight = -4138            # xlMedium on all 4 edges
$rngB28.Borders.Item(10).LineStyle = -4142  # xlEdgeRight -> none

# C28 and D28: blank cells with a full medium box border, default font.
$rngC28 = $ws.Range("C28")
$rngC28.Borders.Weight = -4138

$rngD28 = $ws.Range("D28")
$rngD28.Borders.Weight = -4138

# E28: RMM_C value for 2022 = 57.7, formatted like the rest of column E
# (Trebuchet MS font, centered/top/wrap, 0.0 number format), with a
# medium border on top/right/bottom (no left border).
$rngE28 = $ws.Range("E28")
$rngE28.Value = 57.7
$rngE28.NumberFormat = "0.0"
$rngE28.Font.Name = "Trebuchet MS"
$rngE28.HorizontalAlignment = -4108   # xlCenter
$rngE28.VerticalAlignment = -4160     # xlTop
$rngE28.WrapText = $true
$rngE28.Borders.Weight = -4138
$rngE28.Borders.Item(7).LineStyle = -4142   # xlEdgeLeft -> none

# ---------------------------------------------------------------------
# 3) Update the view: active cell moves to E1, first visible row to A5.
# ---------------------------------------------------------------------
$ws.Range("E1").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1

$wb.Application.Calculate()
